$d = $word.ActiveDocument

# Edit 1: 'There must be at least one sequence within your da'
$full0 = "There must be at least one sequence within your database."
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute($full0, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { Write-Output "EDIT 1: NOT FOUND" }
$paraStart = $r.Start
$offset = 0
$p0_0 = "There must be at "
$segStart = $paraStart + $offset
$segEnd = $segStart + $p0_0.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p0_0.Length
$p0_1 = "least one sequence"
$segStart = $paraStart + $offset
$segEnd = $segStart + $p0_1.Length
$seg = $d.Range($segStart, $segEnd)
$seg.Font.HighlightColorIndex = 7
$offset = $offset + $p0_1.Length
$p0_2 = " within your database."
$segStart = $paraStart + $offset
$segEnd = $segStart + $p0_2.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p0_2.Length

# Edit 2: 'Write the SQL statements that create the indexes y'
$full1 = "Write the SQL statements that create the indexes you specified in objective 1. You must have some in your project."
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute($full1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { Write-Output "EDIT 2: NOT FOUND" }
$paraStart = $r.Start
$offset = 0
$p1_0 = "Write the SQL statements that "
$segStart = $paraStart + $offset
$segEnd = $segStart + $p1_0.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p1_0.Length
$p1_1 = "create the indexes you specified in objective 1. You must have some in your project."
$segStart = $paraStart + $offset
$segEnd = $segStart + $p1_1.Length
$seg = $d.Range($segStart, $segEnd)
$seg.Font.HighlightColorIndex = 7
$offset = $offset + $p1_1.Length

# Edit 3: 'Write SQL statements to INSERT at least 4 or 5 rec'
$full2 = "Write SQL statements to INSERT at least 4 or 5 records into each of your tables. You may want more in some tables. You want to have enough data that allows you to test queries and Transact SQL code that will be developed later. For example you may want more than 4 or 5 records in an intersecting table that exists to resolve a many to many relationship between two tables. "
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute($full2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { Write-Output "EDIT 3: NOT FOUND" }
$paraStart = $r.Start
$offset = 0
$p2_0 = "Write SQL "
$segStart = $paraStart + $offset
$segEnd = $segStart + $p2_0.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p2_0.Length
$p2_1 = "statements to INSERT at least 4 or 5 records"
$segStart = $paraStart + $offset
$segEnd = $segStart + $p2_1.Length
$seg = $d.Range($segStart, $segEnd)
$seg.Font.HighlightColorIndex = 7
$offset = $offset + $p2_1.Length
$p2_2 = " into each of your tables. You may want more in some tables. You want to have enough data that allows you to test queries and Transact SQL code that will be developed later. For example you may want more than 4 or 5 records in an intersecting table that exists to resolve a many to many relationship between two tables. "
$segStart = $paraStart + $offset
$segEnd = $segStart + $p2_2.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p2_2.Length

# Edit 4: 'Create a total of at least six procedures and func'
$full3 = "Create a total of at least six procedures and functions. At least two of these six program units must be functions and at least three must be procedures."
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute($full3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { Write-Output "EDIT 4: NOT FOUND" }
$paraStart = $r.Start
$offset = 0
$p3_0 = "Create a total of at "
$segStart = $paraStart + $offset
$segEnd = $segStart + $p3_0.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p3_0.Length
$p3_1 = "least six procedures and functions"
$segStart = $paraStart + $offset
$segEnd = $segStart + $p3_1.Length
$seg = $d.Range($segStart, $segEnd)
$seg.Font.HighlightColorIndex = 7
$offset = $offset + $p3_1.Length
$p3_2 = ". At least "
$segStart = $paraStart + $offset
$segEnd = $segStart + $p3_2.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p3_2.Length
$p3_3 = "two of these six program units must be functions"
$segStart = $paraStart + $offset
$segEnd = $segStart + $p3_3.Length
$seg = $d.Range($segStart, $segEnd)
$seg.Font.HighlightColorIndex = 7
$offset = $offset + $p3_3.Length
$p3_4 = " and at "
$segStart = $paraStart + $offset
$segEnd = $segStart + $p3_4.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p3_4.Length
$p3_5 = "least three must be procedures"
$segStart = $paraStart + $offset
$segEnd = $segStart + $p3_5.Length
$seg = $d.Range($segStart, $segEnd)
$seg.Font.HighlightColorIndex = 7
$offset = $offset + $p3_5.Length
$p3_6 = "."
$segStart = $paraStart + $offset
$segEnd = $segStart + $p3_6.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p3_6.Length

# Edit 5: 'At least three out of the six program units must c'
$full4 = "At least three out of the six program units must contain EXCEPTION handling. For instance a ‘try…catch’ blocks."
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute($full4, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { Write-Output "EDIT 5: NOT FOUND" }
$paraStart = $r.Start
$offset = 0
$p4_0 = "At least "
$segStart = $paraStart + $offset
$segEnd = $segStart + $p4_0.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p4_0.Length
$p4_1 = "three out of the six program units must contain EXCEPTION handling"
$segStart = $paraStart + $offset
$segEnd = $segStart + $p4_1.Length
$seg = $d.Range($segStart, $segEnd)
$seg.Font.HighlightColorIndex = 7
$offset = $offset + $p4_1.Length
$p4_2 = ". For instance a ‘try…catch’ blocks."
$segStart = $paraStart + $offset
$segEnd = $segStart + $p4_2.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p4_2.Length

# Edit 6: 'In addition to creating the six program units, you'
$full5 = "In addition to creating the six program units, you need to create at least three triggers. One of the triggers should fire on an INSERT, one should fire on an UPDATE and the third should fire on a DELETE. For example you could create triggers that correspond to the INSERT, UPDATE and DELETE procedures mentioned above. The triggers could add information to a log table that is used to track changes to one or more tables."
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute($full5, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { Write-Output "EDIT 6: NOT FOUND" }
$paraStart = $r.Start
$offset = 0
$p5_0 = "In addition to creating the six program units, you need to create at "
$segStart = $paraStart + $offset
$segEnd = $segStart + $p5_0.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p5_0.Length
$p5_1 = "least three triggers"
$segStart = $paraStart + $offset
$segEnd = $segStart + $p5_1.Length
$seg = $d.Range($segStart, $segEnd)
$seg.Font.HighlightColorIndex = 7
$offset = $offset + $p5_1.Length
$p5_2 = ". "
$segStart = $paraStart + $offset
$segEnd = $segStart + $p5_2.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p5_2.Length
$p5_3 = "One of the triggers should fire on an INSERT"
$segStart = $paraStart + $offset
$segEnd = $segStart + $p5_3.Length
$seg = $d.Range($segStart, $segEnd)
$seg.Font.HighlightColorIndex = 7
$offset = $offset + $p5_3.Length
$p5_4 = ", "
$segStart = $paraStart + $offset
$segEnd = $segStart + $p5_4.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p5_4.Length
$p5_5 = "one should fire on an UPDATE"
$segStart = $paraStart + $offset
$segEnd = $segStart + $p5_5.Length
$seg = $d.Range($segStart, $segEnd)
$seg.Font.HighlightColorIndex = 7
$offset = $offset + $p5_5.Length
$p5_6 = " and the "
$segStart = $paraStart + $offset
$segEnd = $segStart + $p5_6.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p5_6.Length
$p5_7 = "third should fire on a DELETE"
$segStart = $paraStart + $offset
$segEnd = $segStart + $p5_7.Length
$seg = $d.Range($segStart, $segEnd)
$seg.Font.HighlightColorIndex = 7
$offset = $offset + $p5_7.Length
$p5_8 = ". For example you could create triggers that correspond to the INSERT, UPDATE and DELETE procedures mentioned above. The triggers could add information to a log table that is used to track changes to one or more tables."
$segStart = $paraStart + $offset
$segEnd = $segStart + $p5_8.Length
$seg = $d.Range($segStart, $segEnd)
$offset = $offset + $p5_8.Length

Write-Output "ALL EDITS DONE"